$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) values for rows 2-30, column G
$newValues = @{
    2  = 5
    3  = 7
    4  = 3
    5  = 3
    6  = 2
    7  = 5
    8  = 8
    9  = 3
    10 = 1
    11 = 1
    12 = 5
    13 = 7
    14 = 5
    15 = 5
    16 = 4
    17 = 3
    18 = 5
    19 = 1
    20 = 1
    21 = 4
    22 = 4
    23 = 7
    24 = 2
    25 = 4
    26 = 2
    27 = 5
    28 = 5
    29 = 2
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
